$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Price" column (D) cells that hold numeric-looking text stay as text
# (they are stored as inline strings, not numbers, in the source workbook).
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D14", "D16", "D17", "D18", "D20", "D21", "D23", "D25", "D26", "D27", "D28", "D30", "D31", "D38", "D39", "D40", "D41", "D43", "D45", "D47", "D48", "D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated values scraped by the GitHub Actions cron job
$ws.Range('D2').Value = '68.385.88'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '2.643.70'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '597.78'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').Value = '154.37'
$ws.Range('E6').Value = '  +1.24%  '
$ws.Range('E8').Value = '  +0.71%  '
$ws.Range('D9').Value = '2.642.75'
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('E10').Value = '  +8.58%  '
$ws.Range('E11').Value = '  -0.53%  '
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('E13').Value = '  +2.33%  '
$ws.Range('D14').Value = '28.17'
$ws.Range('E14').Value = '  +2.42%  '
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').Value = '3.123.94'
$ws.Range('E16').Value = '  +1.16%  '
$ws.Range('D17').Value = '68.363.20'
$ws.Range('E17').Value = '  +1.09%  '
$ws.Range('D18').Value = '2.639.30'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('E19').Value = '  +2.33%  '
$ws.Range('D20').Value = '364.70'
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('D21').Value = '7.50'
$ws.Range('E21').Value = '  +13.60%  '
$ws.Range('E22').Value = '  +3.82%  '
$ws.Range('D23').Value = '4.89'
$ws.Range('E23').Value = '  +2.08%  '
$ws.Range('E24').Value = '  +0.86%  '
$ws.Range('D25').Value = '74.04'
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').Value = '9.81'
$ws.Range('E27').Value = '  -0.44%  '
$ws.Range('D28').Value = '0.0000106'
$ws.Range('E28').Value = '  +2.46%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('D31').Value = '573.28'
$ws.Range('E31').Value = '  -2.87%  '
$ws.Range('E32').Value = '  +4.76%  '
$ws.Range('E34').Value = '  +1.45%  '
$ws.Range('E35').Value = '  +4.47%  '
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('E37').Value = '  +5.47%  '
$ws.Range('D38').Value = '160.35'
$ws.Range('E38').Value = '  +1.00%  '
$ws.Range('D39').Value = '19.37'
$ws.Range('E39').Value = '  +1.17%  '
$ws.Range('D40').Value = '1.91'
$ws.Range('E40').Value = '  +1.19%  '
$ws.Range('D41').Value = '0.373'
$ws.Range('E41').Value = '  +1.78%  '
$ws.Range('E42').Value = '  +2.69%  '
$ws.Range('D43').Value = '0.0₆0342'
$ws.Range('E43').Value = '  +12.65%  '
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').Value = '17.72'
$ws.Range('E45').Value = '  +3.75%  '
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').Value = '40.43'
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').Value = '157.41'
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('E49').Value = '  +2.77%  '
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('D51').Value = '21.92'
$ws.Range('E51').Value = '  +2.93%  '

Write-Host "Updated cryptos list"